$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Correct the footprint on row 19 (Cout) from c0806 to c0805
$ws.Range("F19").Value = "c0805"

# Add a warning note in H10
$ws.Range("H10").Value = "NOG NIET ALLE CAPS EN WEERSTANDEN STAAN IN BOM!!!"

# Rename reference designators U2 -> IC2 and U1 -> IC1 on rows 3 and 2
$ws.Range("A3").Value = "IC2"
$ws.Range("A2").Value = "IC1"

# Give the new empty row below the table the same (hyperlink) look as the
# Link column above it
$ws.Range("E21").Style = "Hyperlink"
$ws.Range("E21").Interior.ColorIndex = -4142

# Select A3 to match the saved selection state
$ws.Range("A3").Select()
